$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.1
$ws.Range("H2").Value = 2.63
$ws.Range("I2").Value = 2.75
$ws.Range("J2").Value = 4.33
$ws.Range("L2").Value = 3.75
$ws.Range("M2").Value = 1.18
$ws.Range("N2").Value = 4.5
$ws.Range("U2").Value = 6.6
$ws.Range("V2").Value = 1.1
$ws.Range("AC2").Value = 5.5
$ws.Range("AD2").Value = 13
$ws.Range("AF2").Value = 34
$ws.Range("AH2").Value = 51
$ws.Range("AI2").Value = 4.33
$ws.Range("AN2").Value = 11
$ws.Range("AO2").Value = 13
$ws.Range("AP2").Value = 29

# Row 3
$ws.Range("G3").Value = 3.2
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 2.05
$ws.Range("J3").Value = 3.75
$ws.Range("K3").Value = 2.3
$ws.Range("L3").Value = 2.63
$ws.Range("M3").Value = 1.03
$ws.Range("N3").Value = 15
$ws.Range("S3").Value = 1.67
$ws.Range("T3").Value = 2.15
$ws.Range("AC3").Value = 13
$ws.Range("AD3").Value = 19
$ws.Range("AJ3").Value = 7
$ws.Range("AM3").Value = 10
$ws.Range("AP3").Value = 19
$ws.Range("AQ3").Value = 15
$ws.Range("AR3").Value = 21
$ws.Range("AS3").Value = 126

# Row 5
$ws.Range("G5").Value = 2.8
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 2.35
$ws.Range("J5").Value = 3.35
$ws.Range("K5").Value = 2.07
$ws.Range("L5").Value = 2.92
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.25
$ws.Range("S5").Value = 1.75
$ws.Range("T5").Value = 1.87
$ws.Range("W5").Value = 2.7
$ws.Range("X5").Value = 1.35
$ws.Range("Z5").Value = 2.52
$ws.Range("AA5").Value = 1.6
$ws.Range("AB5").Value = 2.07
$ws.Range("AC5").Value = 9.75
$ws.Range("AD5").Value = 15
$ws.Range("AE5").Value = 10
$ws.Range("AF5").Value = 35
$ws.Range("AG5").Value = 23
$ws.Range("AI5").Value = 10.75
$ws.Range("AJ5").Value = 6.4
$ws.Range("AL5").Value = 50
$ws.Range("AM5").Value = 9.25
$ws.Range("AO5").Value = 9
$ws.Range("AP5").Value = 25
$ws.Range("AQ5").Value = 18
$ws.Range("AR5").Value = 24
$ws.Range("AS5").Value = 350

# Row 6
$ws.Range("G6").Value = 1.83
$ws.Range("H6").Value = 3.6
$ws.Range("I6").Value = 3.7
$ws.Range("J6").Value = 2.4
$ws.Range("L6").Value = 4
$ws.Range("AA6").Value = 1.67
$ws.Range("AB6").Value = 2.1
$ws.Range("AC6").Value = 8.5
$ws.Range("AD6").Value = 9.5
$ws.Range("AE6").Value = 9
$ws.Range("AG6").Value = 15
$ws.Range("AJ6").Value = 7
$ws.Range("AK6").Value = 13
$ws.Range("AN6").Value = 21
$ws.Range("AO6").Value = 13
$ws.Range("AQ6").Value = 29

# Row 7
$ws.Range("G7").Value = 2.6
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 2.55
$ws.Range("J7").Value = 3.25
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 3.2
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 8
$ws.Range("S7").Value = 2.05
$ws.Range("T7").Value = 1.75
$ws.Range("W7").Value = 3.5
$ws.Range("X7").Value = 1.29
$ws.Range("AC7").Value = 8.5
$ws.Range("AD7").Value = 13
$ws.Range("AF7").Value = 26
$ws.Range("AG7").Value = 23
$ws.Range("AI7").Value = 9
$ws.Range("AJ7").Value = 6.5
$ws.Range("AM7").Value = 8
$ws.Range("AN7").Value = 12
$ws.Range("AQ7").Value = 21
$ws.Range("AS7").Value = 700

# Row 8
$ws.Range("G8").Value = 1.75
$ws.Range("H8").Value = 3.7
$ws.Range("I8").Value = 3.8
$ws.Range("J8").Value = 2.38
$ws.Range("L8").Value = 4.33
$ws.Range("AE8").Value = 9
$ws.Range("AF8").Value = 15
$ws.Range("AH8").Value = 26
$ws.Range("AM8").Value = 12
$ws.Range("AO8").Value = 13
$ws.Range("AR8").Value = 34
$ws.Range("AS8").Value = 600
